$d = $word.ActiveDocument

# --- 1. "Livestock Health Management and Regulation" -> "Office of the Chief Veterinarian"
#        (only the second occurrence, in the footer address block, paragraph 38)
$p = $d.Paragraphs(38)
$p.Range.Find.Execute("Livestock Health Management and Regulation", $true, $false, $false, $false, $false, $true, 1, $false, "Office of the Chief Veterinarian", 2)

# --- 2. "1767 Angus Campbell Road                     Toll: 1 877-877-2474" -> "1767 Angus Campbell Road"
$p = $d.Paragraphs(39)
$p.Range.Find.Execute("1767 Angus Campbell Road                     Toll: 1 877-877-2474", $true, $false, $false, $false, $false, $true, 1, $false, "1767 Angus Campbell Road", 2)

# --- 3. "Abbotsford, B.C.   V3G 2M3                      Tel:  (778) 666-0560  Fax:  (604) 556-3015" -> "Abbotsford, BC   V3G 2M3"
$p = $d.Paragraphs(40)
$p.Range.Find.Execute("Abbotsford, B.C.   V3G 2M3                      Tel:  (778) 666-0560  Fax:  (604) 556-3015", $true, $false, $false, $false, $false, $true, 1, $false, "Abbotsford, BC   V3G 2M3", 2)

# --- 4. "Minister of Finance" gets its own direct font-size override (10.5pt / half-points 21)
$p = $d.Paragraphs(35)
$rng = $p.Range
$rng.Find.Execute("Minister of Finance", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Size = 10.5
